$wb = $excel.ActiveWorkbook

# --- RunManager sheet: add a new test row ---
# ("1" is entered with a leading quote-prefix so it is stored as text, matching
# the existing Priority/Execution Times columns, instead of being auto-typed
# as a number)
$wsRun = $wb.Worksheets.Item("RunManager")
$wsRun.Cells.Item(4, 1).Value = "verifyBooksPageTitle"
$wsRun.Cells.Item(4, 2).Value = "verifyBooksPageTitle"
$wsRun.Cells.Item(4, 3).Value = "yes"
$wsRun.Cells.Item(4, 4).Value = "'1"
$wsRun.Cells.Item(4, 5).Value = "'1"

# --- TestData sheet: add Author/Category-ish "Title" column, flip existing
#     Execution Flag values to "no", and add a row for the new test ---
$wsData = $wb.Worksheets.Item("TestData")

# New header for column G
$wsData.Cells.Item(1, 7).Value = "Title"

# Existing rows 2-6: Execution Flag (col B) all become "no", and col G gets "'"
# (doubled quote escapes the COM quote-prefix so the literal apostrophe is stored)
for ($r = 2; $r -le 6; $r++) {
    $wsData.Cells.Item($r, 2).Value = "no"
    $wsData.Cells.Item($r, 7).Value = "''"
}

# New row 7 for verifyBooksPageTitle test
$wsData.Cells.Item(7, 1).Value = "verifyBooksPageTitle"
$wsData.Cells.Item(7, 2).Value = "yes"
$wsData.Cells.Item(7, 3).Value = "chrome"
$wsData.Cells.Item(7, 4).Value = "''"
$wsData.Cells.Item(7, 5).Value = "''"
$wsData.Cells.Item(7, 6).Value = "''"
$wsData.Cells.Item(7, 7).Value = "Book"
